# Auto-generated Excel COM-interop script
# Updates the cryptos price/volume table (and re-ranks two swapped coin pairs)
# to match the target snapshot described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal TEXT (never let Excel auto-convert
# number-looking strings like "1.000" or "26.499.72" into numeric values),
# while leaving the cell's style/formatting untouched.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.499.72'
Set-TextValue $ws.Range('E2') '  -1.34%  '
Set-TextValue $ws.Range('D3') '1.850.61'
Set-TextValue $ws.Range('E3') '  -1.46%  '
Set-TextValue $ws.Range('D4') '1.000'
Set-TextValue $ws.Range('E4') '  +0.31%  '
Set-TextValue $ws.Range('D5') '260.68'
Set-TextValue $ws.Range('E5') '  -8.12%  '
Set-TextValue $ws.Range('D6') '1.001'
Set-TextValue $ws.Range('E6') '  +0.15%  '
Set-TextValue $ws.Range('D7') '0.5167'
Set-TextValue $ws.Range('E7') '  -0.48%  '
Set-TextValue $ws.Range('D8') '0.3249'
Set-TextValue $ws.Range('E8') '  -8.19%  '
Set-TextValue $ws.Range('D9') '0.06774'
Set-TextValue $ws.Range('E9') '  -4.78%  '
Set-TextValue $ws.Range('D10') '18.95'
Set-TextValue $ws.Range('E10') '  -6.66%  '
Set-TextValue $ws.Range('D11') '0.7723'
Set-TextValue $ws.Range('E11') '  -6.22%  '
Set-TextValue $ws.Range('B12') 'WrappedEther'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D12') '1.899.43'
Set-TextValue $ws.Range('E12') '  +1.20%  '
Set-TextValue $ws.Range('B13') 'TRON'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D13') '0.07720'
Set-TextValue $ws.Range('E13') '  -0.58%  '
Set-TextValue $ws.Range('D14') '88.67'
Set-TextValue $ws.Range('E14') '  -1.45%  '
Set-TextValue $ws.Range('D15') '5.037'
Set-TextValue $ws.Range('E15') '  -2.78%  '
Set-TextValue $ws.Range('D16') '1.000'
Set-TextValue $ws.Range('E16') '  +0.40%  '
Set-TextValue $ws.Range('D17') '14.11'
Set-TextValue $ws.Range('E17') '  -2.65%  '
Set-TextValue $ws.Range('E18') '  +0.18%  '
Set-TextValue $ws.Range('D19') '0.000007919'
Set-TextValue $ws.Range('E19') '  -3.16%  '
Set-TextValue $ws.Range('D20') '26.534.69'
Set-TextValue $ws.Range('E20') '  -1.15%  '
Set-TextValue $ws.Range('D21') '2.113.97'
Set-TextValue $ws.Range('E21') '  +0.62%  '
Set-TextValue $ws.Range('D22') '4.528'
Set-TextValue $ws.Range('E22') '  -5.50%  '
Set-TextValue $ws.Range('D23') '9.542'
Set-TextValue $ws.Range('E23') '  -6.41%  '
Set-TextValue $ws.Range('D24') '5.923'
Set-TextValue $ws.Range('E24') '  -5.14%  '
Set-TextValue $ws.Range('D25') '2.352'
Set-TextValue $ws.Range('E25') '  -3.49%  '
Set-TextValue $ws.Range('D26') '144.58'
Set-TextValue $ws.Range('E26') '  -0.76%  '
Set-TextValue $ws.Range('D27') '1.653'
Set-TextValue $ws.Range('E27') '  -0.90%  '
Set-TextValue $ws.Range('D28') '16.97'
Set-TextValue $ws.Range('E28') '  -2.78%  '
Set-TextValue $ws.Range('D29') '111.23'
Set-TextValue $ws.Range('E29') '  -0.32%  '
Set-TextValue $ws.Range('D30') '4.218'
Set-TextValue $ws.Range('E30') '  -4.77%  '
Set-TextValue $ws.Range('D31') '4.176'
Set-TextValue $ws.Range('E31') '  -4.49%  '
Set-TextValue $ws.Range('D32') '0.08745'
Set-TextValue $ws.Range('E32') '  -1.37%  '
Set-TextValue $ws.Range('E33') '  -2.37%  '
Set-TextValue $ws.Range('D34') '1.134'
Set-TextValue $ws.Range('E34') '  -4.29%  '
Set-TextValue $ws.Range('D35') '2.844'
Set-TextValue $ws.Range('E35') '  -0.75%  '
Set-TextValue $ws.Range('D36') '0.6879'
Set-TextValue $ws.Range('E36') '  -8.29%  '
Set-TextValue $ws.Range('D37') '3.120'
Set-TextValue $ws.Range('E37') '  -5.22%  '
Set-TextValue $ws.Range('D38') '0.01794'
Set-TextValue $ws.Range('E38') '  -4.95%  '
Set-TextValue $ws.Range('D39') '2.211'
Set-TextValue $ws.Range('E39') '  -9.07%  '
Set-TextValue $ws.Range('D40') '0.4903'
Set-TextValue $ws.Range('E40') '  -8.14%  '
Set-TextValue $ws.Range('D41') '112.99'
Set-TextValue $ws.Range('E41') '  -3.11%  '
Set-TextValue $ws.Range('D42') '0.9009'
Set-TextValue $ws.Range('E42') '  -7.77%  '
Set-TextValue $ws.Range('D43') '6.145'
Set-TextValue $ws.Range('E43') '  -2.76%  '
Set-TextValue $ws.Range('D44') '1.000'
Set-TextValue $ws.Range('E44') '  +0.19%  '
Set-TextValue $ws.Range('D45') '7.791'
Set-TextValue $ws.Range('E45') '  -5.22%  '
Set-TextValue $ws.Range('D46') '0.4225'
Set-TextValue $ws.Range('E46') '  -9.06%  '
Set-TextValue $ws.Range('B47') 'Algorand'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D47') '0.1258'
Set-TextValue $ws.Range('E47') '  -8.50%  '
Set-TextValue $ws.Range('B48') 'EnergySwap'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D48') '9.113'
Set-TextValue $ws.Range('E48') '  -4.20%  '
Set-TextValue $ws.Range('D49') '0.05891'
Set-TextValue $ws.Range('E49') '  -0.75%  '
Set-TextValue $ws.Range('D50') '35.27'
Set-TextValue $ws.Range('E50') '  -4.19%  '
Set-TextValue $ws.Range('D51') '59.35'
Set-TextValue $ws.Range('E51') '  -4.31%  '
